$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Adam11"
$ws.Range("C2").Value2 = "Itga4"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.07360033333333334
$ws.Range("H2").Value2 = 0.220801
$ws.Range("I2").Value2 = 0.02873450582079328
$ws.Range("J2").Value2 = 0.02873450582079327
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 26.04517333333333
$ws.Range("N2").Value2 = 78.13552
$ws.Range("O2").Value2 = 0.9210237118384171
$ws.Range("P2").Value2 = 0.921023711838417
$ws.Range("Q2").Value2 = 1.916933439057778
$ws.Range("R2").Value2 = 17.25240095152
$ws.Range("S2").Value2 = 0.02646516120890963
$ws.Range("T2").Value2 = 0.02646516120890962

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Adam11"
$ws.Range("C3").Value2 = "Itga4"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.07360033333333334
$ws.Range("H3").Value2 = 0.220801
$ws.Range("I3").Value2 = 0.02873450582079328
$ws.Range("J3").Value2 = 0.02873450582079327
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.3302223333333333
$ws.Range("N3").Value2 = 0.990667
$ws.Range("O3").Value2 = 0.01167750336256582
$ws.Range("P3").Value2 = 0.01167750336256582
$ws.Range("Q3").Value2 = 0.02430447380744444
$ws.Range("R3").Value2 = 0.218740264267
$ws.Range("S3").Value2 = 0.0003355472883439807
$ws.Range("T3").Value2 = 0.0003355472883439807

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Adam11"
$ws.Range("C4").Value2 = "Itga4"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.07360033333333334
$ws.Range("H4").Value2 = 0.220801
$ws.Range("I4").Value2 = 0.02873450582079328
$ws.Range("J4").Value2 = 0.02873450582079327
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 1.903109
$ws.Range("N4").Value2 = 5.709327
$ws.Range("O4").Value2 = 0.06729878479901708
$ws.Range("P4").Value2 = 0.06729878479901708
$ws.Range("Q4").Value2 = 0.1400694567696667
$ws.Range("R4").Value2 = 1.260625110927
$ws.Range("S4").Value2 = 0.00193379732353967
$ws.Range("T4").Value2 = 0.00193379732353967

$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Adam11"
$ws.Range("C5").Value2 = "Itga4"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 0.4213883333333333
$ws.Range("H5").Value2 = 1.264165
$ws.Range("I5").Value2 = 0.1645153624799848
$ws.Range("J5").Value2 = 0.1645153624799848
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 26.04517333333333
$ws.Range("N5").Value2 = 78.13552
$ws.Range("O5").Value2 = 0.9210237118384171
$ws.Range("P5").Value2 = 0.921023711838417
$ws.Range("Q5").Value2 = 10.97513218231111
$ws.Range("R5").Value2 = 98.7761896408
$ws.Range("S5").Value2 = 0.1515225498057583
$ws.Range("T5").Value2 = 0.1515225498057583

$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Adam11"
$ws.Range("C6").Value2 = "Itga4"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 0.4213883333333333
$ws.Range("H6").Value2 = 1.264165
$ws.Range("I6").Value2 = 0.1645153624799848
$ws.Range("J6").Value2 = 0.1645153624799848
$ws.Range("K6").Value2 = 2
$ws.Range("L6").Value2 = 0.6666666666666666
$ws.Range("M6").Value2 = 0.3302223333333333
$ws.Range("N6").Value2 = 0.990667
$ws.Range("O6").Value2 = 0.01167750336256582
$ws.Range("P6").Value2 = 0.01167750336256582
$ws.Range("Q6").Value2 = 0.1391518386727778
$ws.Range("R6").Value2 = 1.252366548055
$ws.Range("S6").Value2 = 0.001921128698553758
$ws.Range("T6").Value2 = 0.001921128698553758

$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Adam11"
$ws.Range("C7").Value2 = "Itga4"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 0.4213883333333333
$ws.Range("H7").Value2 = 1.264165
$ws.Range("I7").Value2 = 0.1645153624799848
$ws.Range("J7").Value2 = 0.1645153624799848
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 1.903109
$ws.Range("N7").Value2 = 5.709327
$ws.Range("O7").Value2 = 0.06729878479901708
$ws.Range("P7").Value2 = 0.06729878479901708
$ws.Range("Q7").Value2 = 0.8019479296616666
$ws.Range("R7").Value2 = 7.217531366955
$ws.Range("S7").Value2 = 0.01107168397567279
$ws.Range("T7").Value2 = 0.01107168397567279

$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Adam11"
$ws.Range("C8").Value2 = "Itga4"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 2.066403333333334
$ws.Range("H8").Value2 = 6.199210000000001
$ws.Range("I8").Value2 = 0.8067501316992219
$ws.Range("J8").Value2 = 0.8067501316992219
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 26.04517333333333
$ws.Range("N8").Value2 = 78.13552
$ws.Range("O8").Value2 = 0.9210237118384171
$ws.Range("P8").Value2 = 0.921023711838417
$ws.Range("Q8").Value2 = 53.81983299324445
$ws.Range("R8").Value2 = 484.3784969392
$ws.Range("S8").Value2 = 0.7430360008237492
$ws.Range("T8").Value2 = 0.7430360008237492

$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Adam11"
$ws.Range("C9").Value2 = "Itga4"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 2.066403333333334
$ws.Range("H9").Value2 = 6.199210000000001
$ws.Range("I9").Value2 = 0.8067501316992219
$ws.Range("J9").Value2 = 0.8067501316992219
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.3302223333333333
$ws.Range("N9").Value2 = 0.990667
$ws.Range("O9").Value2 = 0.01167750336256582
$ws.Range("P9").Value2 = 0.01167750336256582
$ws.Range("Q9").Value2 = 0.6823725303411112
$ws.Range("R9").Value2 = 6.14135277307
$ws.Range("S9").Value2 = 0.009420827375668085
$ws.Range("T9").Value2 = 0.009420827375668085

$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Adam11"
$ws.Range("C10").Value2 = "Itga4"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 2.066403333333334
$ws.Range("H10").Value2 = 6.199210000000001
$ws.Range("I10").Value2 = 0.8067501316992219
$ws.Range("J10").Value2 = 0.8067501316992219
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 1.903109
$ws.Range("N10").Value2 = 5.709327
$ws.Range("O10").Value2 = 0.06729878479901708
$ws.Range("P10").Value2 = 0.06729878479901708
$ws.Range("Q10").Value2 = 3.932590781296667
$ws.Range("R10").Value2 = 35.39331703167
$ws.Range("S10").Value2 = 0.05429330349980462
$ws.Range("T10").Value2 = 0.05429330349980462
